# Inserts a new weekly data row for "Haba" (Vega Central Mapocho de Santiago)
# at worksheet row 254, pushing the former rows 254-275 down to 255-276.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 254; existing rows 254-275 shift to 255-276.
$ws.Rows.Item(254).Insert()

# Populate the newly inserted row 254 with the new record.
$ws.Cells.Item(254, 1).Value = 9
$ws.Cells.Item(254, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(254, 3).Value = "Metropolitana"
$ws.Cells.Item(254, 4).Value = 44826
$ws.Cells.Item(254, 4).NumberFormat = $ws.Cells.Item(255, 4).NumberFormat
$ws.Cells.Item(254, 5).Value = 13
$ws.Cells.Item(254, 6).Value = 100112026
$ws.Cells.Item(254, 7).Value = "Haba"
$ws.Cells.Item(254, 8).Value = "Sin especificar"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 80
$ws.Cells.Item(254, 11).Value = 10000
$ws.Cells.Item(254, 12).Value = 10000
$ws.Cells.Item(254, 13).Value = 10000
$ws.Cells.Item(254, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(254, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(254, 16).Value = 400
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"
